$wb = $excel.ActiveWorkbook

# Drop the duplicate "test_sheet" tab, keeping "sheet1" (which preserves
# the original sheetId/rId) and renaming it to "test_sheet".
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("test_sheet").Delete()

$ws = $wb.Worksheets.Item("sheet1")
$ws.Name = "test_sheet"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 11
$ws.Range("C2").Value = 21

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "#N/A"
$ws.Range("C3").Value = 22

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 13
$ws.Range("C4").Value = 23
